# "upload superived train 0.2k/0.4k/2k/5k"
#
# 1) Add a 6th data series ("3k"-ish extra mAP numbers) into column G,
#    rows 2-16, alongside the existing map/aps/apm/apl table (A1:E16).
# 2) Replace the old scratch area below (rows 22, 24, 27-29) with a new,
#    properly laid out comparison table in rows 22-26: a header row of
#    epoch labels (1k..15k) plus four method rows ("more ways", "scar",
#    "block", "3 way") of mAP values.
# 3) Reposition/resize the chart object.
# 4) Update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G values for the existing table (rows 2-16) ---
$ws.Range("G2").Value  = 0.277
$ws.Range("G3").Value  = 0.395
$ws.Range("G4").Value  = 0.383
$ws.Range("G5").Value  = 0.479
$ws.Range("G6").Value  = 0.474
$ws.Range("G7").Value  = 0.527
$ws.Range("G8").Value  = 0.527
$ws.Range("G9").Value  = 0.549
$ws.Range("G10").Value = 0.538
$ws.Range("G11").Value = 0.547
$ws.Range("G12").Value = 0.553
$ws.Range("G13").Value = 0.554
$ws.Range("G14").Value = 0.546
$ws.Range("G15").Value = 0.578
$ws.Range("G16").Value = 0.578

# --- Clear the old scratch rows (22, 24, 27-29) so nothing stray remains ---
$ws.Range("A22:P29").ClearContents()

# --- Row 22: header with epoch labels ---
$ws.Range("A22").Value = "cutpaste methods"
$ws.Range("B22").Value = "1k"
$ws.Range("C22").Value = "2k"
$ws.Range("D22").Value = "3k"
$ws.Range("E22").Value = "4k"
$ws.Range("F22").Value = "5k"
$ws.Range("G22").Value = "6k"
$ws.Range("H22").Value = "7k"
$ws.Range("I22").Value = "8k"
$ws.Range("J22").Value = "9k"
$ws.Range("K22").Value = "10k"
$ws.Range("L22").Value = "11k"
$ws.Range("M22").Value = "12k"
$ws.Range("N22").Value = "13k"
$ws.Range("O22").Value = "14k"
$ws.Range("P22").Value = "15k"

# --- Row 23: "more ways" ---
$ws.Range("A23").Value = "more ways"
$ws.Range("B23").Value = 0.302
$ws.Range("C23").Value = 0.317
$ws.Range("D23").Value = 0.42
$ws.Range("E23").Value = 0.481
$ws.Range("F23").Value = 0.484
$ws.Range("G23").Value = 0.519
$ws.Range("H23").Value = 0.494
$ws.Range("I23").Value = 0.518
$ws.Range("J23").Value = 0.558
$ws.Range("K23").Value = 0.556
$ws.Range("L23").Value = 0.587
$ws.Range("M23").Value = 0.585
$ws.Range("N23").Value = 0.607
$ws.Range("O23").Value = 0.602
$ws.Range("P23").Value = 0.597

# --- Row 24: "scar" ---
$ws.Range("A24").Value = "scar"
$ws.Range("B24").Value = 0.27
$ws.Range("C24").Value = 0.282
$ws.Range("D24").Value = 0.302
$ws.Range("E24").Value = 0.352
$ws.Range("F24").Value = 0.353
$ws.Range("G24").Value = 0.383
$ws.Range("H24").Value = 0.402
$ws.Range("I24").Value = 0.417
$ws.Range("J24").Value = 0.42
$ws.Range("K24").Value = 0.46
$ws.Range("L24").Value = 0.446
$ws.Range("M24").Value = 0.492
$ws.Range("N24").Value = 0.55
$ws.Range("O24").Value = 0.554
$ws.Range("P24").Value = 0.553

# --- Row 25: "block" ---
$ws.Range("A25").Value = "block"
$ws.Range("B25").Value = 0.28
$ws.Range("C25").Value = 0.303
$ws.Range("D25").Value = 0.332
$ws.Range("E25").Value = 0.34
$ws.Range("F25").Value = 0.338
$ws.Range("G25").Value = 0.37
$ws.Range("H25").Value = 0.422
$ws.Range("I25").Value = 0.448
$ws.Range("J25").Value = 0.484
$ws.Range("K25").Value = 0.49
$ws.Range("L25").Value = 0.512
$ws.Range("M25").Value = 0.549
$ws.Range("N25").Value = 0.568
$ws.Range("O25").Value = 0.567
$ws.Range("P25").Value = 0.566

# --- Row 26: "3 way" ---
$ws.Range("A26").Value = "3 way"
$ws.Range("B26").Value = 0.277
$ws.Range("C26").Value = 0.312
$ws.Range("D26").Value = 0.34
$ws.Range("E26").Value = 0.342
$ws.Range("F26").Value = 0.35
$ws.Range("G26").Value = 0.42
$ws.Range("H26").Value = 0.467
$ws.Range("I26").Value = 0.47
$ws.Range("J26").Value = 0.488
$ws.Range("K26").Value = 0.52
$ws.Range("L26").Value = 0.544
$ws.Range("M26").Value = 0.55
$ws.Range("N26").Value = 0.561
$ws.Range("O26").Value = 0.562
$ws.Range("P26").Value = 0.56

# --- Reposition the chart (was roughly G6:M21, now roughly H2:O17) ---
$co = $ws.ChartObjects().Item(1)
$co.Left   = 461.9911
$co.Top    = 17.1499
$co.Width  = 385.0625
$co.Height = 216.1071

# --- Update the selected cell to match the saved view state ---
$ws.Range("J28").Select()
